$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data
$ws.Range("D2").Value = '27.730.65'
$ws.Range("E2").Value = '  +3.05%  '
$ws.Range("D3").Value = '1.864.13'
$ws.Range("E3").Value = '  +2.95%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.033'
$ws.Range("E4").Value = '  +2.59%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.16'
$ws.Range("E5").Value = '  +3.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.032'
$ws.Range("E6").Value = '  +2.59%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4415'
$ws.Range("E7").Value = '  +2.85%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3801'
$ws.Range("E8").Value = '  +2.87%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07450'
$ws.Range("E9").Value = '  +2.77%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8840'
$ws.Range("E10").Value = '  +2.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.83'
$ws.Range("E11").Value = '  +2.95%  '
$ws.Range("D12").Value = '1.878.11'
$ws.Range("E12").Value = '  -8.73%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.558'
$ws.Range("E13").Value = '  +3.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.750'
$ws.Range("E14").Value = '  +1.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07199'
$ws.Range("E15").Value = '  +4.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.89'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.038'
$ws.Range("E17").Value = '  +3.10%  '
$ws.Range("E18").Value = '  +1.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.032'
$ws.Range("E19").Value = '  +2.61%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.53'
$ws.Range("E20").Value = '  +2.35%  '
$ws.Range("D21").Value = '27.770.01'
$ws.Range("E21").Value = '  +3.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.310'
$ws.Range("E22").Value = '  +2.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.43'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '158.45'
$ws.Range("E24").Value = '  +2.86%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.932'
$ws.Range("E25").Value = '  +2.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.85'
$ws.Range("E26").Value = '  +2.89%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.335'
$ws.Range("E27").Value = '  +1.76%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.988'
$ws.Range("E28").Value = '  +4.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '117.71'
$ws.Range("E29").Value = '  +2.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09093'
$ws.Range("E30").Value = '  +1.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.214'
$ws.Range("E31").Value = '  +5.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7683'
$ws.Range("E32").Value = '  +3.65%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.570'
$ws.Range("E33").Value = '  +3.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.943'
$ws.Range("E34").Value = '  +4.93%  '
$ws.Range("E35").Value = '  +2.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.157'
$ws.Range("E36").Value = '  +2.81%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01991'
$ws.Range("E37").Value = '  +3.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05345'
$ws.Range("E38").Value = '  +2.26%  '
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5206'
$ws.Range("E39").Value = '  +2.34%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.833'
$ws.Range("E40").Value = '  +3.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1693'
$ws.Range("E41").Value = '  +2.73%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.852'
$ws.Range("E42").Value = '  +6.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.732'
$ws.Range("E43").Value = '  +5.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '109.64'
$ws.Range("E44").Value = '  +2.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.57'
$ws.Range("E45").Value = '  +1.63%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.732'
$ws.Range("E46").Value = '  +5.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4690'
$ws.Range("E47").Value = '  +2.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06420'
$ws.Range("E48").Value = '  +2.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.870'
$ws.Range("E49").Value = '  +3.28%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '39.77'
$ws.Range("E50").Value = '  +4.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9359'
$ws.Range("E51").Value = '  +2.71%  '
